$d = $word.ActiveDocument

$replacements = @(
    @{old = "2025-06-29 Sunday"; new = "2025-06-30 Monday"},
    @{old = "710÷4="; new = "601÷4="},
    @{old = "460÷8="; new = "835÷4="},
    @{old = "494÷6="; new = "240÷7="},
    @{old = "492÷3="; new = "464÷5="},
    @{old = "283÷9="; new = "305÷6="},
    @{old = "998÷4="; new = "722÷2="},
    @{old = "745÷5="; new = "196÷6="},
    @{old = "363÷3="; new = "882÷6="},
    @{old = "608÷2="; new = "110÷8="},
    @{old = "119÷4="; new = "334÷5="},
    @{old = "861÷7="; new = "518÷2="},
    @{old = "550÷4="; new = "929÷5="},
    @{old = "260÷7="; new = "294÷6="},
    @{old = "213÷9="; new = "931÷6="},
    @{old = "988÷6="; new = "868÷9="},
    @{old = "896÷2="; new = "745÷3="},
    @{old = "352÷5="; new = "803÷5="},
    @{old = "970÷5="; new = "658÷8="},
    @{old = "188÷5="; new = "309÷8="},
    @{old = "768÷2="; new = "778÷4="},
    @{old = "437÷7="; new = "504÷2="},
    @{old = "684÷5="; new = "280÷4="},
    @{old = "898÷4="; new = "982÷6="},
    @{old = "341÷8="; new = "128÷4="},
    @{old = "880÷9="; new = "415÷9="}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $false, $false, $false, $false, `
                       $true, 1, $false, $r.new, 2)
}

$d.Save()
